# issue #5: add legislator_id, name, date into dataframe
#
# The 股票 (stock) sheet is the only sheet that gains new columns: three
# trailing columns (date, legislator_name, legislator_id) are appended
# after the existing "total" column, both as a new header (row 1) and as
# data (row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# ---- header row (row 1) -------------------------------------------------
# Copy the formatting of an existing header cell (bold font, thin border,
# centered/top aligned) onto the three new header cells, then fill in the
# text so they match the look of B1:G1.
$ws.Range("B1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# ---- data row (row 2) ----------------------------------------------------
# legislator_name / legislator_id are plain text/number, so they pick up
# the surrounding (unformatted) data-row look with no extra work.
$ws.Range("I2").Value = "許智傑"
$ws.Range("J2").Value = 1750

# date needs special handling: "2013-12-11" typed into a General cell
# would be auto-parsed into a date serial number, but the source data
# keeps it as a literal text string. Force text formatting first, assign
# the value, then drop the formatting override again so the cell's style
# doesn't keep the temporary "text" number format attached.
$dateCell = $ws.Range("H2")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2013-12-11"
$dateCell.ClearFormats()
